# Collapse the split "<id>" / "<ident>" / "</id>" runs back into a single
# run for each of the p112r_1 .. p112r_4 tagged-id blocks (the fig_p112r_1
# block is left untouched, matching the source diff).
$d = $word.ActiveDocument

$ids = @("p112r_1", "p112r_2", "p112r_3", "p112r_4")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2) | Out-Null
}
